$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Перевірив" (reviewer) line: fill in the reviewer's name after the
#    underlined blank, and relocate the "_GoBack" bookmark to sit right
#    after the newly typed name (Word drops the old "_GoBack" bookmark
#    automatically once a new one with the same name is added).
# ---------------------------------------------------------------------
$reviewRng = $d.Content
$reviewRng.Find.ClearFormatting()
$blankAfterReviewed = "Перевірив                      "
$found = $reviewRng.Find.Execute($blankAfterReviewed, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: could not locate the reviewer blank"
}
$reviewRng.Collapse(0)
$reviewRng.InsertAfter("Вєчерковська Анастасія Сергіївна")
$reviewRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $reviewRng)

# ---------------------------------------------------------------------
# 2) Mark the five picture runs (C++ / python screenshots) as
#    "do not spell-check" (<w:noProof/>) to match how Word tags runs
#    that hold inline drawings once it has re-laid them out.
# ---------------------------------------------------------------------
$pictureIndexes = @(1, 2, 8, 9, 10)
foreach ($i in $pictureIndexes) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = 1
}

# ---------------------------------------------------------------------
# 3) Tail sentence: "...надавалося" + "." used to be split by the
#    (now-relocated) "_GoBack" bookmark; merge them back into a single
#    run reading "...надавалося.".
# ---------------------------------------------------------------------
$tailText = $d.Content.Text
$tailIdx = $tailText.IndexOf("надавалося")
$tailLen = "надавалося.".Length
$tailRange = $d.Range($tailIdx, $tailIdx + $tailLen)
$tailRange.Text = "Q"
$afterPlaceholder = $d.Content.Text
$placeholderIdx = $afterPlaceholder.IndexOf("Q")
$placeholderRange = $d.Range($placeholderIdx, $placeholderIdx + 1)
$placeholderRange.Text = "надавалося."

Write-Host "edit complete"
